# Weekly update: insert three new daily records (2021-09-14) above the
# existing history, pushing all prior rows down by three.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows starting at row 3 (old rows 3:21 shift to 6:24).
$ws.Rows("3:5").Insert()

# New data for the three inserted rows.
$newRows = @(
    @(3, "Femacal de La Calera", "Coquimbo", 44453, 5, "Fruta", 100107, "Otros", 100107002, "Chirimoya", "Cultivar IV Región", "Especial", 45, 30000, 30000, 30000, "`$/bandeja 10 kilos", "Provincia del Elquí", 3000, 10),
    @(3, "Femacal de La Calera", "Coquimbo", 44453, 5, "Fruta", 100107, "Otros", 100107002, "Chirimoya", "Cultivar IV Región", "Primera", 47, 27000, 27000, 27000, "`$/bandeja 10 kilos", "Provincia del Elquí", 2700, 10),
    @(3, "Femacal de La Calera", "Coquimbo", 44453, 5, "Fruta", 100107, "Otros", 100107002, "Chirimoya", "Cultivar IV Región", "Segunda", 40, 25000, 25000, 25000, "`$/bandeja 10 kilos", "Provincia del Elquí", 2500, 10)
)

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = 3 + $i
    $rowData = $newRows[$i]
    for ($c = 1; $c -le $rowData.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowData[$c - 1]
    }
}
